$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Drop the last two slides (id 264 "The grass is always greener on
#    the other side" and id 265, the chart/table slide) together with
#    their notes pages, which PowerPoint removes automatically.
# ---------------------------------------------------------------------
while ($p.Slides.Count -gt 17) {
    $p.Slides.Item($p.Slides.Count).Delete()
}

# ---------------------------------------------------------------------
# 2. Re-stamp the cached "datetimeFigureOut" placeholder text (the
#    auto date field) from 2/1/2024 to 2/4/2024 everywhere it is
#    cached: the slide master, every slide layout, and the notes
#    master.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "2/1/2024") {
                $shp.TextFrame.TextRange.Text = "2/4/2024"
            }
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

Update-DatePlaceholder $p.NotesMaster.Shapes
